$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last week's rows (312:315) down to new rows (316:319)
# so the previous week's data is preserved further down the sheet.
$ws.Range("A312:R315").Copy()
$ws.Range("A316:R319").Insert()

# Update rows 312:315 in place with the new week's values (date 44595)
$ws.Range("D312").Value = 44595
$ws.Range("I312").Value = "Primera"
$ws.Range("J312").Value = 3000
$ws.Range("K312").Value = 500
$ws.Range("L312").Value = 500
$ws.Range("M312").Value = 500
$ws.Range("P312").Value = 500

$ws.Range("D313").Value = 44595
$ws.Range("I313").Value = "Segunda"
$ws.Range("J313").Value = 3000
$ws.Range("K313").Value = 300
$ws.Range("L313").Value = 300
$ws.Range("M313").Value = 300
$ws.Range("P313").Value = 300

$ws.Range("D314").Value = 44595
$ws.Range("I314").Value = "Primera"
$ws.Range("J314").Value = 3000
$ws.Range("K314").Value = 500
$ws.Range("L314").Value = 500
$ws.Range("M314").Value = 500
$ws.Range("P314").Value = 500

$ws.Range("D315").Value = 44595
$ws.Range("I315").Value = "Segunda"
$ws.Range("J315").Value = 3000
$ws.Range("K315").Value = 300
$ws.Range("L315").Value = 300
$ws.Range("M315").Value = 300
$ws.Range("P315").Value = 300
